{"js": "// Office.js (Word JavaScript API) edit script.\n// Updates the date heading and every arithmetic-expression cell in the\n// single table, in document order, replacing old text with new text.\nconst replacements = [\n  [\"2022-12-17 Saturday\", \"2022-12-18 Sunday\"],\n  [\"4+34=\", \"96-57=\"],\n  [\"53+41=\", \"81-22=\"],\n  [\"72+10=\", \"27+22=\"],\n  [\"79-11=\", \"86-75=\"],\n  [\"46-31=\", \"2+38=\"],\n  [\"14+2=\", \"40+26=\"],\n  [\"37-34=\", \"12+84=\"],\n  [\"39+17=\", \"94+1=\"],\n  [\"43+44=\", \"26-7=\"],\n  [\"93-34=\", \"37+41=\"],\n  [\"13-7=\", \"8+87=\"],\n  [\"10+41=\", \"21-18=\"],\n  [\"18-3=\", \"64-30=\"],\n  [\"67-13=\", \"53-20=\"],\n  [\"83-43=\", \"61-5=\"],\n  [\"17+73=\", \"34+17=\"],\n  [\"59-59=\", \"73-58=\"],\n  [\"48+18=\", \"99-64=\"],\n  [\"48-2=\", \"95-53=\"],\n  [\"99-72=\", \"17-12=\"],\n  [\"33+56=\", \"74-63=\"],\n  [\"11+74=\", \"58+6=\"],\n  [\"60-51=\", \"85-75=\"],\n  [\"22+68=\", \"82-9=\"],\n  [\"59-12=\", \"14+34=\"],\n  [\"68-19=\", \"34-2=\"],\n  [\"1+96=\", \"57-1=\"],\n  [\"87-32=\", \"93-87=\"],\n  [\"88-79=\", \"6-0=\"],\n  [\"8+9=\", \"91-47=\"],\n  [\"68-50=\", \"37+10=\"],\n  [\"86-26=\", \"83+13=\"],\n  [\"77-49=\", \"19+0=\"],\n  [\"30+66=\", \"90-61=\"],\n  [\"60+8=\", \"44+36=\"],\n  [\"93-77=\", \"15+28=\"],\n  [\"37+54=\", \"42+40=\"],\n  [\"66-31=\", \"34+52=\"],\n  [\"21+51=\", \"72-7=\"],\n  [\"55-29=\", \"61-51=\"],\n  [\"95-13=\", \"54-7=\"],\n  [\"68-45=\", \"67+28=\"],\n  [\"94-69=\", \"94-60=\"],\n  [\"14+0=\", \"81-37=\"],\n  [\"13+15=\", \"44-34=\"],\n  [\"84-62=\", \"46+4=\"],\n  [\"29+57=\", \"16+47=\"],\n  [\"19+79=\", \"94-87=\"],\n  [\"20+43=\", \"68-53=\"],\n  [\"22+52=\", \"67-21=\"],\n  [\"59-0=\", \"30+17=\"],\n  [\"34+0=\", \"72-5=\"],\n  [\"41-8=\", \"94-54=\"],\n  [\"57-45=\", \"48-16=\"],\n  [\"94-72=\", \"6+92=\"],\n  [\"27+26=\", \"92-50=\"],\n  [\"5+50=\", \"17+53=\"],\n  [\"67+7=\", \"52+4=\"],\n  [\"48-31=\", \"88-19=\"],\n  [\"61+36=\", \"74-36=\"],\n  [\"23-19=\", \"23+57=\"],\n  [\"3+19=\", \"61-13=\"],\n  [\"78+5=\", \"15+13=\"],\n  [\"61-45=\", \"38+10=\"],\n  [\"1+13=\", \"74+22=\"],\n  [\"58-47=\", \"88-84=\"],\n  [\"54-9=\", \"85-37=\"],\n  [\"44-28=\", \"28+63=\"],\n  [\"72-62=\", \"49-48=\"],\n  [\"77-32=\", \"37+29=\"],\n  [\"90-35=\", \"16+29=\"],\n  [\"73-67=\", \"31-10=\"],\n  [\"67-16=\", \"66+13=\"],\n  [\"19+60=\", \"65+26=\"],\n  [\"1+50=\", \"94-21=\"],\n  [\"97-86=\", \"13+13=\"],\n  [\"72-36=\", \"25+22=\"],\n  [\"44+8=\", \"39+9=\"],\n  [\"20+22=\", \"9+31=\"],\n  [\"40-3=\", \"61-22=\"],\n  [\"86-70=\", \"61-8=\"],\n  [\"56-9=\", \"64-46=\"],\n  [\"25+67=\", \"78-68=\"],\n  [\"43+36=\", \"90-3=\"],\n  [\"12+9=\", \"32+35=\"],\n  [\"28-18=\", \"78-76=\"],\n  [\"4+22=\", \"28+27=\"],\n  [\"18+39=\", \"70-49=\"],\n  [\"66-17=\", \"24+56=\"],\n  [\"0+26=\", \"87-75=\"],\n  [\"15+6=\", \"15+63=\"],\n  [\"47-19=\", \"73-26=\"],\n  [\"54-19=\", \"67+4=\"],\n  [\"95-75=\", \"78-24=\"],\n  [\"9+87=\", \"38+5=\"],\n  [\"39-23=\", \"57-12=\"],\n  [\"6+93=\", \"45-10=\"],\n  [\"62+1=\", \"27+36=\"],\n  [\"91-34=\", \"62-1=\"],\n  [\"70-22=\", \"81-68=\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== replacements.length) {\n  throw new Error(\n    `Expected ${replacements.length} paragraphs, found ${paragraphs.items.length}`\n  );\n}\n\n// Load text for every paragraph up front so we can verify alignment\n// before mutating anything.\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const para = paragraphs.items[i];\n  if (para.text !== oldText) {\n    throw new Error(\n      `Paragraph ${i} text mismatch: expected \"${oldText}\", found \"${para.text}\"`\n    );\n  }\n  para.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Updates the date heading and every arithmetic-expression cell in the\n# single table. We walk $d.Paragraphs in document order (this includes the\n# empty paragraph-mark placeholders Word inserts at the end of each table\n# row) and rewrite only the 101 non-empty paragraphs' text in place by\n# paragraph position. This avoids Find/Replace entirely, which is important\n# here because several *new* values contain an *old* value as a literal\n# substring (e.g. old '0+26=' is a substring of new '40+26='); a plain\n# Find-based replace pass would let an earlier replacement's result get\n# re-matched (and corrupted) by a later search.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2022-12-17 Saturday\", \"2022-12-18 Sunday\"),\n    @(\"4+34=\", \"96-57=\"),\n    @(\"53+41=\", \"81-22=\"),\n    @(\"72+10=\", \"27+22=\"),\n    @(\"79-11=\", \"86-75=\"),\n    @(\"46-31=\", \"2+38=\"),\n    @(\"14+2=\", \"40+26=\"),\n    @(\"37-34=\", \"12+84=\"),\n    @(\"39+17=\", \"94+1=\"),\n    @(\"43+44=\", \"26-7=\"),\n    @(\"93-34=\", \"37+41=\"),\n    @(\"13-7=\", \"8+87=\"),\n    @(\"10+41=\", \"21-18=\"),\n    @(\"18-3=\", \"64-30=\"),\n    @(\"67-13=\", \"53-20=\"),\n    @(\"83-43=\", \"61-5=\"),\n    @(\"17+73=\", \"34+17=\"),\n    @(\"59-59=\", \"73-58=\"),\n    @(\"48+18=\", \"99-64=\"),\n    @(\"48-2=\", \"95-53=\"),\n    @(\"99-72=\", \"17-12=\"),\n    @(\"33+56=\", \"74-63=\"),\n    @(\"11+74=\", \"58+6=\"),\n    @(\"60-51=\", \"85-75=\"),\n    @(\"22+68=\", \"82-9=\"),\n    @(\"59-12=\", \"14+34=\"),\n    @(\"68-19=\", \"34-2=\"),\n    @(\"1+96=\", \"57-1=\"),\n    @(\"87-32=\", \"93-87=\"),\n    @(\"88-79=\", \"6-0=\"),\n    @(\"8+9=\", \"91-47=\"),\n    @(\"68-50=\", \"37+10=\"),\n    @(\"86-26=\", \"83+13=\"),\n    @(\"77-49=\", \"19+0=\"),\n    @(\"30+66=\", \"90-61=\"),\n    @(\"60+8=\", \"44+36=\"),\n    @(\"93-77=\", \"15+28=\"),\n    @(\"37+54=\", \"42+40=\"),\n    @(\"66-31=\", \"34+52=\"),\n    @(\"21+51=\", \"72-7=\"),\n    @(\"55-29=\", \"61-51=\"),\n    @(\"95-13=\", \"54-7=\"),\n    @(\"68-45=\", \"67+28=\"),\n    @(\"94-69=\", \"94-60=\"),\n    @(\"14+0=\", \"81-37=\"),\n    @(\"13+15=\", \"44-34=\"),\n    @(\"84-62=\", \"46+4=\"),\n    @(\"29+57=\", \"16+47=\"),\n    @(\"19+79=\", \"94-87=\"),\n    @(\"20+43=\", \"68-53=\"),\n    @(\"22+52=\", \"67-21=\"),\n    @(\"59-0=\", \"30+17=\"),\n    @(\"34+0=\", \"72-5=\"),\n    @(\"41-8=\", \"94-54=\"),\n    @(\"57-45=\", \"48-16=\"),\n    @(\"94-72=\", \"6+92=\"),\n    @(\"27+26=\", \"92-50=\"),\n    @(\"5+50=\", \"17+53=\"),\n    @(\"67+7=\", \"52+4=\"),\n    @(\"48-31=\", \"88-19=\"),\n    @(\"61+36=\", \"74-36=\"),\n    @(\"23-19=\", \"23+57=\"),\n    @(\"3+19=\", \"61-13=\"),\n    @(\"78+5=\", \"15+13=\"),\n    @(\"61-45=\", \"38+10=\"),\n    @(\"1+13=\", \"74+22=\"),\n    @(\"58-47=\", \"88-84=\"),\n    @(\"54-9=\", \"85-37=\"),\n    @(\"44-28=\", \"28+63=\"),\n    @(\"72-62=\", \"49-48=\"),\n    @(\"77-32=\", \"37+29=\"),\n    @(\"90-35=\", \"16+29=\"),\n    @(\"73-67=\", \"31-10=\"),\n    @(\"67-16=\", \"66+13=\"),\n    @(\"19+60=\", \"65+26=\"),\n    @(\"1+50=\", \"94-21=\"),\n    @(\"97-86=\", \"13+13=\"),\n    @(\"72-36=\", \"25+22=\"),\n    @(\"44+8=\", \"39+9=\"),\n    @(\"20+22=\", \"9+31=\"),\n    @(\"40-3=\", \"61-22=\"),\n    @(\"86-70=\", \"61-8=\"),\n    @(\"56-9=\", \"64-46=\"),\n    @(\"25+67=\", \"78-68=\"),\n    @(\"43+36=\", \"90-3=\"),\n    @(\"12+9=\", \"32+35=\"),\n    @(\"28-18=\", \"78-76=\"),\n    @(\"4+22=\", \"28+27=\"),\n    @(\"18+39=\", \"70-49=\"),\n    @(\"66-17=\", \"24+56=\"),\n    @(\"0+26=\", \"87-75=\"),\n    @(\"15+6=\", \"15+63=\"),\n    @(\"47-19=\", \"73-26=\"),\n    @(\"54-19=\", \"67+4=\"),\n    @(\"95-75=\", \"78-24=\"),\n    @(\"9+87=\", \"38+5=\"),\n    @(\"39-23=\", \"57-12=\"),\n    @(\"6+93=\", \"45-10=\"),\n    @(\"62+1=\", \"27+36=\"),\n    @(\"91-34=\", \"62-1=\"),\n    @(\"70-22=\", \"81-68=\"),\n)\n\n$n = $d.Paragraphs.Count\n$idx = 0\nfor ($i = 1; $i -le $n; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n    $clean = $r.Text -replace \"[\\r\\x07]\", \"\"\n    if ($clean -eq \"\") {\n        continue\n    }\n\n    if ($idx -ge $replacements.Length) {\n        throw \"More non-empty paragraphs than expected replacements\"\n    }\n\n    $oldText = $replacements[$idx][0]\n    $newText = $replacements[$idx][1]\n    $idx = $idx + 1\n\n    if ($clean -ne $oldText) {\n        throw \"Paragraph $i text mismatch: expected [$oldText], found [$clean]\"\n    }\n\n    # Trim the trailing paragraph mark (and cell mark, if any) from the\n    # range before assigning so only the run text is replaced; formatting\n    # on the existing run (font / size) is preserved.\n    $markLen = $r.End - $r.Start - $clean.Length\n    $textRange = $d.Range($r.Start, $r.End - $markLen)\n    $textRange.Text = $newText\n}\n\nif ($idx -ne $replacements.Length) {\n    throw \"Expected $($replacements.Length) non-empty paragraphs, found $idx\"\n}\n\n"}
